$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells are written as Text so numeric-looking
# strings (e.g. "226.36") are not auto-converted to numbers by Excel,
# matching the inlineStr text cells in the source workbook.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.139.51"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.789.17"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.36"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("E6").Value = "  -0.63%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.94"
$ws.Range("E8").Value = "  -0.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.293"
$ws.Range("E9").Value = "  +1.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0691"
$ws.Range("E10").Value = "  -1.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0946"
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.047.29"
$ws.Range("E12").Value = "  +0.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.17"
$ws.Range("E13").Value = "  +2.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.796.43"
$ws.Range("E14").Value = "  +1.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.076.65"
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.621"
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("E17").Value = "  +1.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.08"
$ws.Range("E18").Value = "  +0.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.03"
$ws.Range("E19").Value = "  +1.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0779"
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("E22").Value = "  +1.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.10"
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.05"
$ws.Range("E24").Value = "  -0.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.44"
$ws.Range("E25").Value = "  +0.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.15"
$ws.Range("E26").Value = "  +1.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.31"
$ws.Range("E28").Value = "  +0.74%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("E31").Value = "  +1.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.66"
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.62"
$ws.Range("E33").Value = "  +3.18%  "
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.452.87"
$ws.Range("E35").Value = "  +4.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0192"
$ws.Range("E38").Value = "  +2.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.04"
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.26"
$ws.Range("E40").Value = "  +3.63%  "
$ws.Range("E41").Value = "  +0.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.920"
$ws.Range("E42").Value = "  +1.33%  "
$ws.Range("E43").Value = "  +0.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.50"
$ws.Range("E44").Value = "  +2.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0510"
$ws.Range("E45").Value = "  +2.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.05"
$ws.Range("E46").Value = "  +4.33%  "
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("E48").Value = "  -1.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.948.51"
$ws.Range("E49").Value = "  +1.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "106.19"
$ws.Range("E50").Value = "  -1.45%  "
$ws.Range("E51").Value = "  -0.06%  "

# Rows 36 and 37 swapped ranking: ImmutableX now ranks above RenderToken
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.648"
$ws.Range("E36").Value = "  -0.55%  "

$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.42"
$ws.Range("E37").Value = "  +8.48%  "
